$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.179.12'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '1.855.26'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.94'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6896'
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07774'
$ws.Range("E8").Value = '  +3.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3046'
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.20'
$ws.Range("E10").Value = '  -2.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08063'
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("D12").Value = '1.843.11'
$ws.Range("E12").Value = '  -2.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7212'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.190'
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.30'
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '29.182.73'
$ws.Range("E16").Value = '  -2.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.738'
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007806'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.26'
$ws.Range("E19").Value = '  +0.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '234.75'
$ws.Range("E20").Value = '  -3.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("D22").Value = '2.111.75'
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.472'
$ws.Range("E24").Value = '  -2.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.10'
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.972'
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("E27").Value = '  -4.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.04'
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.010'
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05205'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.184'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7031'
$ws.Range("E36").Value = '  -2.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.015'
$ws.Range("E37").Value = '  +1.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.673'
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01849'
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.680'
$ws.Range("E40").Value = '  -1.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9434'
$ws.Range("E41").Value = '  +6.35%  '
$ws.Range("D42").Value = '1.093.54'
$ws.Range("E42").Value = '  +4.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.943'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4287'
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.43'
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9997'
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.39'
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("D49").Value = '2.007.87'
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.165'
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.005'
$ws.Range("E51").Value = '  -3.47%  '
